# Update crypto price/volume snapshot values (GitHub Actions data refresh).
# For D-column prices that parse as plain numbers, force text format first so
# Excel doesn't silently coerce the literal string into a numeric value (which
# would lose formatting such as trailing zeros, e.g. "61.40" -> 61.4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.619.82'
$ws.Range('E2').Value = '  +2.42%  '
$ws.Range('D3').Value = '2.097.15'
$ws.Range('E3').Value = '  +3.30%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.13'
$ws.Range('E5').Value = '  +0.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.615'
$ws.Range('E6').Value = '  +1.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '61.29'
$ws.Range('E7').Value = '  +2.41%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +1.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0843'
$ws.Range('E10').Value = '  +1.96%  '
$ws.Range('E11').Value = '  +0.85%  '
$ws.Range('D12').Value = '2.408.87'
$ws.Range('E12').Value = '  +3.53%  '
$ws.Range('E13').Value = '  +1.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.43'
$ws.Range('E14').Value = '  +6.70%  '
$ws.Range('E15').Value = '  +2.15%  '
$ws.Range('E16').Value = '  +5.61%  '
$ws.Range('D17').Value = '2.175.99'
$ws.Range('E17').Value = '  +7.27%  '
$ws.Range('D18').Value = '38.515.19'
$ws.Range('E18').Value = '  +2.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '70.84'
$ws.Range('E19').Value = '  +1.97%  '
$ws.Range('E20').Value = '  +2.14%  '
$ws.Range('D21').Value = '0.0₃0834'
$ws.Range('E21').Value = '  +1.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '227.05'
$ws.Range('E22').Value = '  +1.47%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('E24').Value = '  +1.84%  '
$ws.Range('E25').Value = '  +3.65%  '
$ws.Range('E26').Value = '  +1.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.42'
$ws.Range('E27').Value = '  +0.56%  '
$ws.Range('E28').Value = '  +1.46%  '
$ws.Range('E29').Value = '  +1.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.36'
$ws.Range('E30').Value = '  +8.35%  '
$ws.Range('E31').Value = '  -0.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.32'
$ws.Range('E32').Value = '  +5.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.77'
$ws.Range('E33').Value = '  +6.75%  '
$ws.Range('E34').Value = '  +2.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0605'
$ws.Range('E35').Value = '  +0.28%  '
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.40'
$ws.Range('E37').Value = '  +4.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.55'
$ws.Range('E38').Value = '  +4.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.52'
$ws.Range('E40').Value = '  +2.66%  '
$ws.Range('D41').Value = '1.542.50'
$ws.Range('E41').Value = '  +0.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.93'
$ws.Range('E42').Value = '  +4.63%  '
$ws.Range('E43').Value = '  +2.97%  '
$ws.Range('E44').Value = '  +1.31%  '
$ws.Range('E45').Value = '  +0.40%  '
$ws.Range('E46').Value = '  +2.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.56'
$ws.Range('E47').Value = '  +6.58%  '
$ws.Range('E48').Value = '  +1.54%  '
$ws.Range('E49').Value = '  +3.59%  '
$ws.Range('E50').Value = '  +0.63%  '
$ws.Range('D51').Value = '2.293.79'
$ws.Range('E51').Value = '  +3.42%  '
